$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '258.81'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.50%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.91'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-1.97%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.682'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '2.14%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06002'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.95%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.664'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.52%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8576'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9256'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '0.09%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1390'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-1.43%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.04813'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '33.48%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07013'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.89%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03123'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-3.22%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09133'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.52%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001532'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.42%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006051'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-94.21%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006079'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.26%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.461'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.61%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.164'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.07%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.165'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-1.75%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.13%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.59%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.123'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '7.25%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04238'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.55%'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.19%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-6.25%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001199'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.03%'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '13.56%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03840'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.21%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1115'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '1.15%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.003844'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-38.17%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002418'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '9.96%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01526'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '29.68%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005110'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-6.41%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000749'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.03%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-16.71%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1503'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '16.39%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.03%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.03%'
